$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column D (shifts existing quarterly
# columns D:K to F:M) and pull the formatting for the new columns from
# the columns immediately to their right so number formats/styles match.
$ws.Range("D1:E1").EntireColumn.Insert()
$ws.Range("F7:G102").Copy()
$ws.Range("D7:E102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# New quarter figures (2019-01-31 in D, 2018-10-31 in E) for every
# populated row; rows left out here stay blank like their neighbours.
$newQuarterData = @{
    7 = @(43496, 43404)
    8 = @(23600, 22900)
    9 = @(19100, 18800)
    10 = @(4500, 4100)
    12 = @("NA", "NA")
    13 = @(0, 0)
    14 = @(0, 0)
    15 = @(0, 0)
    17 = @(22900, 22600)
    18 = @(700, 300)
    20 = @(0, 0)
    21 = @(900, 600)
    22 = @(100, 100)
    23 = @(600, 300)
    24 = @(100, 100)
    25 = @(0, 0)
    26 = @(500, 200)
    27 = @(300, 100)
    28 = @(0, 0)
    29 = @(0, 0)
    30 = @(0, 0)
    31 = @(0, 0)
    32 = @(0, 0)
    33 = @(300, 100)
    34 = @(0, 0)
    35 = @(300, 100)
    38 = @(43496, 43404)
    41 = @(4800, 4600)
    42 = @(0, 0)
    43 = @(9400, 10300)
    44 = @(15400, 15300)
    45 = @(600, 600)
    46 = @(30100, 30800)
    47 = @(100, 100)
    48 = @(2300, 2400)
    49 = @(4600, 4600)
    50 = @(0, 0)
    51 = @(0, 0)
    52 = @(1200, 1000)
    53 = @(0, 0)
    54 = @(38300, 38800)
    57 = @(3800, 4800)
    58 = @(5800, 6300)
    59 = @(600, 0)
    60 = @(10200, 11200)
    61 = @(0, 0)
    62 = @(1600, 1700)
    63 = @(0, 0)
    64 = @(0, 0)
    65 = @(0, 0)
    66 = @(13100, 14000)
    68 = @(0, 0)
    69 = @(0, 0)
    70 = @(0, 0)
    71 = @(0, 0)
    72 = @(13700, 13400)
    73 = @(0, 0)
    74 = @(0, 0)
    75 = @(0, 0)
    76 = @(25200, 24900)
    77 = @(0, 0)
    80 = @(43496, 43404)
    81 = @(300, 100)
    83 = @(200, 200)
    84 = @(0, 0)
    85 = @(0, 0)
    86 = @(0, 0)
    87 = @(0, 0)
    88 = @(0, 0)
    89 = @(800, 700)
    91 = @(-100, -100)
    92 = @(0, 0)
    93 = @(0, 0)
    94 = @(-100, 0)
    96 = @(0, 0)
    97 = @(0, 0)
    98 = @(0, 0)
    99 = @(0, 0)
    100 = @(-500, -1200)
    101 = @(0, 0)
    102 = @(200, -500)
}

foreach ($r in $newQuarterData.Keys) {
    $vals = $newQuarterData[$r]
    $ws.Cells.Item($r, 4).Value = $vals[0]
    $ws.Cells.Item($r, 5).Value = $vals[1]
}
